$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 43148.457
$ws.Range("I28").Value = 51299.75
$ws.Range("K28").Value = 51299.75
$ws.Range("M28").Value = -50814.75

$ws.Range("H43").Value = 2015
$ws.Range("J43").Value = 2064.2856
$ws.Range("L43").Value = 2064.2856
$ws.Range("N43").Value = -2202.2856

$ws.Range("H51").Value = 9349.429
$ws.Range("I51").Value = 9315.666999999999
$ws.Range("J51").Value = 9374.75
$ws.Range("K51").Value = 9315.666999999999
$ws.Range("L51").Value = 9374.75
$ws.Range("M51").Value = -8831.666999999999
$ws.Range("N51").Value = -10342.75

$ws.Range("H106").Value = 103397.1
$ws.Range("I106").Value = 3863.125
$ws.Range("J106").Value = 501533
$ws.Range("K106").Value = 3863.125
$ws.Range("L106").Value = 501533
$ws.Range("M106").Value = -3232.125
$ws.Range("N106").Value = -502795

$ws.Range("H111").Value = 13357.777
$ws.Range("J111").Value = 55415.5
$ws.Range("L111").Value = 166246.5
$ws.Range("N111").Value = -172380.5

$ws.Range("H121").Value = 2188.3333
$ws.Range("J121").Value = 2188.3333
$ws.Range("L121").Value = 6564.999899999999
$ws.Range("N121").Value = -10058.9999

$ws.Range("H132").Value = 27031996
$ws.Range("I132").Value = 28575868
$ws.Range("K132").Value = 85727604
$ws.Range("M132").Value = -85725074

$ws.Range("H135").Value = 612.3871
$ws.Range("I135").Value = 519.5714
$ws.Range("J135").Value = 807.3
$ws.Range("K135").Value = 4676.1426
$ws.Range("L135").Value = 7265.7
$ws.Range("M135").Value = -2141.1426
$ws.Range("N135").Value = -12335.7

$ws.Range("H137").Value = 3724.8484
$ws.Range("I137").Value = 2792.5
$ws.Range("J137").Value = 6211.1113
$ws.Range("K137").Value = 8377.5
$ws.Range("L137").Value = 18633.3339
$ws.Range("M137").Value = -5827.5
$ws.Range("N137").Value = -23733.3339

$ws.Range("H138").Value = 3970.1924
$ws.Range("I138").Value = 3992
$ws.Range("J138").Value = 3956.5625
$ws.Range("K138").Value = 11976
$ws.Range("L138").Value = 11869.6875
$ws.Range("M138").Value = -6836
$ws.Range("N138").Value = -22149.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21638.625
$ws.Range("I32").Value = 12952.182
$ws.Range("J32").Value = 40748.8
$ws.Range("K32").Value = 12952.182
$ws.Range("L32").Value = 40748.8
$ws.Range("M32").Value = -12665.182
$ws.Range("N32").Value = -41322.8

$ws.Range("H45").Value = 328506.28
$ws.Range("I45").Value = 398560.22
$ws.Range("K45").Value = 398560.22
$ws.Range("M45").Value = -398183.22

$ws.Range("H61").Value = 3423.4736
$ws.Range("I61").Value = 2708.5881
$ws.Range("K61").Value = 2708.5881
$ws.Range("M61").Value = -2496.5881

$ws.Range("H63").Value = 54549450
$ws.Range("I63").Value = 71430700
$ws.Range("J63").Value = 25007250
$ws.Range("K63").Value = 71430700
$ws.Range("L63").Value = 25007250
$ws.Range("M63").Value = -71430014
$ws.Range("N63").Value = -25008622

$ws.Range("H66").Value = 54549450
$ws.Range("I66").Value = 71430700
$ws.Range("J66").Value = 25007250
$ws.Range("K66").Value = 357153500
$ws.Range("L66").Value = 125036250
$ws.Range("M66").Value = -357150068
$ws.Range("N66").Value = -125043114

$ws.Range("H74").Value = 33337494
$ws.Range("I74").Value = 58826576
$ws.Range("J74").Value = 5614.769
$ws.Range("K74").Value = 58826576
$ws.Range("L74").Value = 5614.769
$ws.Range("M74").Value = -58825702
$ws.Range("N74").Value = -7362.769

$ws.Range("H77").Value = 33337494
$ws.Range("I77").Value = 58826576
$ws.Range("J77").Value = 5614.769
$ws.Range("K77").Value = 294132880
$ws.Range("L77").Value = 28073.845
$ws.Range("M77").Value = -294128512
$ws.Range("N77").Value = -36809.845

$ws.Range("H97").Value = 961.45
$ws.Range("I97").Value = 846.1111
$ws.Range("J97").Value = 1999.5
$ws.Range("K97").Value = 846.1111
$ws.Range("L97").Value = 1999.5
$ws.Range("M97").Value = -350.1111
$ws.Range("N97").Value = -2991.5

$ws.Range("H122").Value = 3107.5386
$ws.Range("J122").Value = 5127.533
$ws.Range("L122").Value = 15382.599
$ws.Range("N122").Value = -20282.599

$ws.Range("H136").Value = 3423.4736
$ws.Range("I136").Value = 2708.5881
$ws.Range("K136").Value = 8125.7643
$ws.Range("M136").Value = -5575.7643

$ws.Range("H137").Value = 48836.668
$ws.Range("J137").Value = 48836.668
$ws.Range("L137").Value = 48836.668
$ws.Range("N137").Value = -59036.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 100002880
$ws.Range("I107").Value = 3500
$ws.Range("K107").Value = 3500
$ws.Range("M107").Value = -1580

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1348.3158
$ws.Range("I105").Value = 1333.125
$ws.Range("K105").Value = 1333.125
$ws.Range("M105").Value = 413.875

$ws.Range("H107").Value = 885.8570999999999
$ws.Range("I107").Value = 928.9048
$ws.Range("K107").Value = 928.9048
$ws.Range("M107").Value = 991.0952

$ws.Range("H122").Value = 3527.3333
$ws.Range("I122").Value = 3514
$ws.Range("J122").Value = 3714
$ws.Range("K122").Value = 10542
$ws.Range("L122").Value = 11142
$ws.Range("M122").Value = -8092
$ws.Range("N122").Value = -16042

$ws.Range("H132").Value = 219169.94
$ws.Range("I132").Value = 1354.1177
$ws.Range("J132").Value = 836314.75
$ws.Range("K132").Value = 4062.3531
$ws.Range("L132").Value = 2508944.25
$ws.Range("M132").Value = -1532.3531
$ws.Range("N132").Value = -2514004.25

$ws.Range("H134").Value = 3108.2334
$ws.Range("I134").Value = 2545
$ws.Range("K134").Value = 7635
$ws.Range("M134").Value = -5100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 57580030
$ws.Range("I4").Value = 111456160
$ws.Range("J4").Value = 43525390
$ws.Range("K4").Value = 334368480
$ws.Range("L4").Value = 130576170
$ws.Range("M4").Value = -334368368
$ws.Range("N4").Value = -130576394

$ws.Range("H20").Value = 24.75
$ws.Range("J20").Value = 24.75
$ws.Range("L20").Value = 74.25
$ws.Range("N20").Value = -528.25

$ws.Range("H24").Value = 92.14286
$ws.Range("I24").Value = 39
$ws.Range("J24").Value = 225
$ws.Range("K24").Value = 117
$ws.Range("L24").Value = 675
$ws.Range("M24").Value = 113
$ws.Range("N24").Value = -1135

$ws.Range("H25").Value = 616.8333
$ws.Range("J25").Value = 800.3333
$ws.Range("L25").Value = 2400.9999
$ws.Range("N25").Value = -2738.9999

$ws.Range("H30").Value = 616.8333
$ws.Range("J30").Value = 800.3333
$ws.Range("L30").Value = 2400.9999
$ws.Range("N30").Value = -2604.9999

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H113").Value = 1389.4062
$ws.Range("J113").Value = 1538.3889
$ws.Range("L113").Value = 4615.1667
$ws.Range("N113").Value = -8955.1667

$ws.Range("H131").Value = 13389.9
$ws.Range("J131").Value = 16946.867
$ws.Range("L131").Value = 50840.601
$ws.Range("N131").Value = -60920.601

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 594.7143
$ws.Range("I107").Value = 591.6667
$ws.Range("K107").Value = 591.6667
$ws.Range("M107").Value = 1328.3333

$ws.Range("H122").Value = 181471.03
$ws.Range("I122").Value = 207800.11
$ws.Range("J122").Value = 3749.75
$ws.Range("K122").Value = 623400.33
$ws.Range("L122").Value = 11249.25
$ws.Range("M122").Value = -620950.33
$ws.Range("N122").Value = -16149.25

$ws.Range("H126").Value = 8078.9062
$ws.Range("I126").Value = 11501.923
$ws.Range("K126").Value = 34505.769
$ws.Range("M126").Value = -32035.769

$ws.Range("H132").Value = 4052.6099
$ws.Range("I132").Value = 3336.4412
$ws.Range("J132").Value = 7531.143
$ws.Range("K132").Value = 10009.3236
$ws.Range("L132").Value = 22593.429
$ws.Range("M132").Value = -7479.3236
$ws.Range("N132").Value = -27653.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 224928.5
$ws.Range("J125").Value = 224928.5
$ws.Range("L125").Value = 224928.5
$ws.Range("N125").Value = -234768.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2186.5
$ws.Range("I122").Value = 1709.7333
$ws.Range("K122").Value = 5129.199900000001
$ws.Range("M122").Value = -2679.199900000001

$ws.Range("H136").Value = 7816.8945
$ws.Range("I136").Value = 11109.059
$ws.Range("J136").Value = 5151.8096
$ws.Range("K136").Value = 33327.177
$ws.Range("L136").Value = 15455.4288
$ws.Range("M136").Value = -30777.177
$ws.Range("N136").Value = -20555.4288
